# Edit described by the target diff:
#   1) Three tables (slides 14, 15, 16) switch from the custom "Table_0"
#      table style ({3FDCD588-D7A1-44CA-B689-335E37F1B29F}) to the built-in
#      table style {BD84568D-77F6-423D-9266-399925387E4F}.
#   2) The presentation's theme colour palette ("Red Violet"/Integral,
#      physically stored in ppt/theme/theme2.xml, which is the theme that
#      the single slide master in this deck actually uses) is replaced by
#      the "Office" colour palette that ppt/theme/theme1.xml already has.

$p = $ppt.ActivePresentation

# --- 1) Re-style the three tables -------------------------------------
$newStyleId = "{BD84568D-77F6-423D-9266-399925387E4F}"
foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    $tableShape = $slide.Shapes.Item(1)
    $tableShape.Table.ApplyStyle($newStyleId)
}

# --- 2) Swap the live theme's colour scheme over to the "Office" colours
$colorScheme = $p.Slides.Item(1).ThemeColorScheme
$colorScheme.Item(1).RGB  = 0        # dk1      000000
$colorScheme.Item(2).RGB  = 16777215 # lt1      FFFFFF
$colorScheme.Item(3).RGB  = 6968388  # dk2      44546A
$colorScheme.Item(4).RGB  = 15132391 # lt2      E7E6E6
$colorScheme.Item(5).RGB  = 13998939 # accent1  5B9BD5
$colorScheme.Item(6).RGB  = 3243501  # accent2  ED7D31
$colorScheme.Item(7).RGB  = 10855845 # accent3  A5A5A5
$colorScheme.Item(8).RGB  = 49407    # accent4  FFC000
$colorScheme.Item(9).RGB  = 12874308 # accent5  4472C4
$colorScheme.Item(10).RGB = 4697456  # accent6  70AD47
$colorScheme.Item(11).RGB = 12673797 # hlink    0563C1
$colorScheme.Item(12).RGB = 7491477  # folHlink 954F72
